$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Update row 5 odds values as per the diff
$ws.Range("G5").Value = 3
$ws.Range("H5").Value = 3.25
$ws.Range("I5").Value = 2.3
$ws.Range("J5").Value = 3.5
$ws.Range("K5").Value = 2.2
$ws.Range("M5").Value = 1.05
$ws.Range("N5").Value = 11
$ws.Range("O5").Value = 1.25
$ws.Range("P5").Value = 3.75
$ws.Range("Q5").Value = 1.83
$ws.Range("R5").Value = 1.98
$ws.Range("S5").Value = 1.36
$ws.Range("T5").Value = 3
$ws.Range("W5").Value = 11
$ws.Range("Y5").Value = 11
$ws.Range("Z5").Value = 29
$ws.Range("AC5").Value = 11
$ws.Range("AM5").Value = 26
$ws.Range("AO5").Value = 15
$ws.Range("AT5").Value = 3
$ws.Range("AW5").Value = 4.5
